$d = $word.ActiveDocument

# --- Change 1: contact line paragraph (add github link + shrink font to 10pt) ---
$p2 = $d.Paragraphs(2)
$p2.Range.Font.Size = 10
$p2.Range.Font.SizeBi = 10

$d.Content.Find.Execute(
    "www.danielshamany.com | dshamany@gmail.com",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "www.danielshamany.com | github.com/dshamany | dshamany@gmail.com",
    2
) | Out-Null

# --- Change 2: Freelance Software Engineer bullet list ---
$d.Content.Find.Execute(
    "- Wrote back-end Python/Django CRUD functionality",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Utilized NodeJS, ExpressJS, PassportJS for most projects, currently using ReactJS/React Native^l- Wrote back-end Python/Django/AWS S3 with CRUD functionality",
    2
) | Out-Null

# --- Change 3: CRM backend bullet ---
$d.Content.Find.Execute(
    "- Developed a CRM backend utilizing AWS, SuiteCRM, MySQL, and RHEL7",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Developed a CRM backend utilizing AWS EC2, SuiteCRM, NodeJS, ExpressJS, REST APIs, ",
    2
) | Out-Null
